# Apply updated TPM-based NATMI ligand-receptor scores to Sheet1.
# Rewrites data rows 2-21 (Tgfb2 -> Tgfbr3 pairs) with recalculated values and
# adds a new "MuSCs" sending-cluster block (rows 17-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Tgfb2"
$ws.Cells.Item(2,3).Value2 = "Tgfbr3"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 2.767474
$ws.Cells.Item(2,8).Value2 = 8.302422
$ws.Cells.Item(2,9).Value2 = 0.1192484523516842
$ws.Cells.Item(2,10).Value2 = 0.1192484523516842
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 41.607769
$ws.Cells.Item(2,14).Value2 = 124.823307
$ws.Cells.Item(2,15).Value2 = 0.4674897261536314
$ws.Cells.Item(2,16).Value2 = 0.4674897261536314
$ws.Cells.Item(2,17).Value2 = 115.148418905506
$ws.Cells.Item(2,18).Value2 = 1036.335770149554
$ws.Cells.Item(2,19).Value2 = 0.0557474263341332
$ws.Cells.Item(2,20).Value2 = 0.0557474263341332

# Row 3
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Tgfb2"
$ws.Cells.Item(3,3).Value2 = "Tgfbr3"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 2.767474
$ws.Cells.Item(3,8).Value2 = 8.302422
$ws.Cells.Item(3,9).Value2 = 0.1192484523516842
$ws.Cells.Item(3,10).Value2 = 0.1192484523516842
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 40.85529333333333
$ws.Cells.Item(3,14).Value2 = 122.56588
$ws.Cells.Item(3,15).Value2 = 0.4590351838457449
$ws.Cells.Item(3,16).Value2 = 0.4590351838457449
$ws.Cells.Item(3,17).Value2 = 113.0659620623733
$ws.Cells.Item(3,18).Value2 = 1017.59365856136
$ws.Cells.Item(3,19).Value2 = 0.05473923524857589
$ws.Cells.Item(3,20).Value2 = 0.05473923524857589

# Row 4
$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Tgfb2"
$ws.Cells.Item(4,3).Value2 = "Tgfbr3"
$ws.Cells.Item(4,4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 2.767474
$ws.Cells.Item(4,8).Value2 = 8.302422
$ws.Cells.Item(4,9).Value2 = 0.1192484523516842
$ws.Cells.Item(4,10).Value2 = 0.1192484523516842
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 0.2401453333333333
$ws.Cells.Item(4,14).Value2 = 0.720436
$ws.Cells.Item(4,15).Value2 = 0.002698185430636104
$ws.Cells.Item(4,16).Value2 = 0.002698185430636104
$ws.Cells.Item(4,17).Value2 = 0.6645959662213333
$ws.Cells.Item(4,18).Value2 = 5.981363695992
$ws.Cells.Item(4,19).Value2 = 0.0003217544367612178
$ws.Cells.Item(4,20).Value2 = 0.0003217544367612179

# Row 5
$ws.Cells.Item(5,1).Value2 = "ECs"
$ws.Cells.Item(5,2).Value2 = "Tgfb2"
$ws.Cells.Item(5,3).Value2 = "Tgfbr3"
$ws.Cells.Item(5,4).Value2 = "MuSCs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 2.767474
$ws.Cells.Item(5,8).Value2 = 8.302422
$ws.Cells.Item(5,9).Value2 = 0.1192484523516842
$ws.Cells.Item(5,10).Value2 = 0.1192484523516842
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 6.288874666666666
$ws.Cells.Item(5,14).Value2 = 18.866624
$ws.Cells.Item(5,15).Value2 = 0.0706595034147231
$ws.Cells.Item(5,16).Value2 = 0.0706595034147231
$ws.Cells.Item(5,17).Value2 = 17.40429712925867
$ws.Cells.Item(5,18).Value2 = 156.638674163328
$ws.Cells.Item(5,19).Value2 = 0.008426036426144272
$ws.Cells.Item(5,20).Value2 = 0.008426036426144272

# Row 6
$ws.Cells.Item(6,1).Value2 = "ECs"
$ws.Cells.Item(6,2).Value2 = "Tgfb2"
$ws.Cells.Item(6,3).Value2 = "Tgfbr3"
$ws.Cells.Item(6,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 2.767474
$ws.Cells.Item(6,8).Value2 = 8.302422
$ws.Cells.Item(6,9).Value2 = 0.1192484523516842
$ws.Cells.Item(6,10).Value2 = 0.1192484523516842
$ws.Cells.Item(6,11).Value2 = 1
$ws.Cells.Item(6,12).Value2 = 0.3333333333333333
$ws.Cells.Item(6,13).Value2 = 0.010449
$ws.Cells.Item(6,14).Value2 = 0.031347
$ws.Cells.Item(6,15).Value2 = 0.0001174011552645203
$ws.Cells.Item(6,16).Value2 = 0.0001174011552645203
$ws.Cells.Item(6,17).Value2 = 0.028917335826
$ws.Cells.Item(6,18).Value2 = 0.260256022434
$ws.Cells.Item(6,19).Value2 = 0.00001399990606959382
$ws.Cells.Item(6,20).Value2 = 0.00001399990606959383

# Row 7
$ws.Cells.Item(7,1).Value2 = "FAPs"
$ws.Cells.Item(7,2).Value2 = "Tgfb2"
$ws.Cells.Item(7,3).Value2 = "Tgfbr3"
$ws.Cells.Item(7,4).Value2 = "ECs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 12.88577866666667
$ws.Cells.Item(7,8).Value2 = 38.657336
$ws.Cells.Item(7,9).Value2 = 0.5552388796954726
$ws.Cells.Item(7,10).Value2 = 0.5552388796954726
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 41.607769
$ws.Cells.Item(7,14).Value2 = 124.823307
$ws.Cells.Item(7,15).Value2 = 0.4674897261536314
$ws.Cells.Item(7,16).Value2 = 0.4674897261536314
$ws.Cells.Item(7,17).Value2 = 536.1485021477946
$ws.Cells.Item(7,18).Value2 = 4825.336519330152
$ws.Cells.Item(7,19).Value2 = 0.2595684718186856
$ws.Cells.Item(7,20).Value2 = 0.2595684718186856

# Row 8
$ws.Cells.Item(8,1).Value2 = "FAPs"
$ws.Cells.Item(8,2).Value2 = "Tgfb2"
$ws.Cells.Item(8,3).Value2 = "Tgfbr3"
$ws.Cells.Item(8,4).Value2 = "FAPs"
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 12.88577866666667
$ws.Cells.Item(8,8).Value2 = 38.657336
$ws.Cells.Item(8,9).Value2 = 0.5552388796954726
$ws.Cells.Item(8,10).Value2 = 0.5552388796954726
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 40.85529333333333
$ws.Cells.Item(8,14).Value2 = 122.56588
$ws.Cells.Item(8,15).Value2 = 0.4590351838457449
$ws.Cells.Item(8,16).Value2 = 0.4590351838457449
$ws.Cells.Item(8,17).Value2 = 526.4522672550755
$ws.Cells.Item(8,18).Value2 = 4738.07040529568
$ws.Cells.Item(8,19).Value2 = 0.2548741812193167
$ws.Cells.Item(8,20).Value2 = 0.2548741812193167

# Row 9
$ws.Cells.Item(9,1).Value2 = "FAPs"
$ws.Cells.Item(9,2).Value2 = "Tgfb2"
$ws.Cells.Item(9,3).Value2 = "Tgfbr3"
$ws.Cells.Item(9,4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 12.88577866666667
$ws.Cells.Item(9,8).Value2 = 38.657336
$ws.Cells.Item(9,9).Value2 = 0.5552388796954726
$ws.Cells.Item(9,10).Value2 = 0.5552388796954726
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 0.2401453333333333
$ws.Cells.Item(9,14).Value2 = 0.720436
$ws.Cells.Item(9,15).Value2 = 0.002698185430636104
$ws.Cells.Item(9,16).Value2 = 0.002698185430636104
$ws.Cells.Item(9,17).Value2 = 3.094459613166222
$ws.Cells.Item(9,18).Value2 = 27.850136518496
$ws.Cells.Item(9,19).Value2 = 0.001498137455717037
$ws.Cells.Item(9,20).Value2 = 0.001498137455717037

# Row 10
$ws.Cells.Item(10,1).Value2 = "FAPs"
$ws.Cells.Item(10,2).Value2 = "Tgfb2"
$ws.Cells.Item(10,3).Value2 = "Tgfbr3"
$ws.Cells.Item(10,4).Value2 = "MuSCs"
$ws.Cells.Item(10,5).Value2 = 3
$ws.Cells.Item(10,6).Value2 = 1
$ws.Cells.Item(10,7).Value2 = 12.88577866666667
$ws.Cells.Item(10,8).Value2 = 38.657336
$ws.Cells.Item(10,9).Value2 = 0.5552388796954726
$ws.Cells.Item(10,10).Value2 = 0.5552388796954726
$ws.Cells.Item(10,11).Value2 = 3
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 6.288874666666666
$ws.Cells.Item(10,14).Value2 = 18.866624
$ws.Cells.Item(10,15).Value2 = 0.0706595034147231
$ws.Cells.Item(10,16).Value2 = 0.0706595034147231
$ws.Cells.Item(10,17).Value2 = 81.03704701707376
$ws.Cells.Item(10,18).Value2 = 729.333423153664
$ws.Cells.Item(10,19).Value2 = 0.03923290351582927
$ws.Cells.Item(10,20).Value2 = 0.03923290351582927

# Row 11
$ws.Cells.Item(11,1).Value2 = "FAPs"
$ws.Cells.Item(11,2).Value2 = "Tgfb2"
$ws.Cells.Item(11,3).Value2 = "Tgfbr3"
$ws.Cells.Item(11,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(11,5).Value2 = 3
$ws.Cells.Item(11,6).Value2 = 1
$ws.Cells.Item(11,7).Value2 = 12.88577866666667
$ws.Cells.Item(11,8).Value2 = 38.657336
$ws.Cells.Item(11,9).Value2 = 0.5552388796954726
$ws.Cells.Item(11,10).Value2 = 0.5552388796954726
$ws.Cells.Item(11,11).Value2 = 1
$ws.Cells.Item(11,12).Value2 = 0.3333333333333333
$ws.Cells.Item(11,13).Value2 = 0.010449
$ws.Cells.Item(11,14).Value2 = 0.031347
$ws.Cells.Item(11,15).Value2 = 0.0001174011552645203
$ws.Cells.Item(11,16).Value2 = 0.0001174011552645203
$ws.Cells.Item(11,17).Value2 = 0.134643501288
$ws.Cells.Item(11,18).Value2 = 1.211791511592
$ws.Cells.Item(11,19).Value2 = 0.00006518568592402649
$ws.Cells.Item(11,20).Value2 = 0.00006518568592402649

# Row 12
$ws.Cells.Item(12,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(12,2).Value2 = "Tgfb2"
$ws.Cells.Item(12,3).Value2 = "Tgfbr3"
$ws.Cells.Item(12,4).Value2 = "ECs"
$ws.Cells.Item(12,5).Value2 = 1
$ws.Cells.Item(12,6).Value2 = 0.3333333333333333
$ws.Cells.Item(12,7).Value2 = 0.03423166666666667
$ws.Cells.Item(12,8).Value2 = 0.102695
$ws.Cells.Item(12,9).Value2 = 0.001475017749550216
$ws.Cells.Item(12,10).Value2 = 0.001475017749550216
$ws.Cells.Item(12,11).Value2 = 3
$ws.Cells.Item(12,12).Value2 = 1
$ws.Cells.Item(12,13).Value2 = 41.607769
$ws.Cells.Item(12,14).Value2 = 124.823307
$ws.Cells.Item(12,15).Value2 = 0.4674897261536314
$ws.Cells.Item(12,16).Value2 = 0.4674897261536314
$ws.Cells.Item(12,17).Value2 = 1.424303279151667
$ws.Cells.Item(12,18).Value2 = 12.818729512365
$ws.Cells.Item(12,19).Value2 = 0.0006895556438089762
$ws.Cells.Item(12,20).Value2 = 0.0006895556438089763

# Row 13
$ws.Cells.Item(13,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(13,2).Value2 = "Tgfb2"
$ws.Cells.Item(13,3).Value2 = "Tgfbr3"
$ws.Cells.Item(13,4).Value2 = "FAPs"
$ws.Cells.Item(13,5).Value2 = 1
$ws.Cells.Item(13,6).Value2 = 0.3333333333333333
$ws.Cells.Item(13,7).Value2 = 0.03423166666666667
$ws.Cells.Item(13,8).Value2 = 0.102695
$ws.Cells.Item(13,9).Value2 = 0.001475017749550216
$ws.Cells.Item(13,10).Value2 = 0.001475017749550216
$ws.Cells.Item(13,11).Value2 = 3
$ws.Cells.Item(13,12).Value2 = 1
$ws.Cells.Item(13,13).Value2 = 40.85529333333333
$ws.Cells.Item(13,14).Value2 = 122.56588
$ws.Cells.Item(13,15).Value2 = 0.4590351838457449
$ws.Cells.Item(13,16).Value2 = 0.4590351838457449
$ws.Cells.Item(13,17).Value2 = 1.398544782955556
$ws.Cells.Item(13,18).Value2 = 12.5869030466
$ws.Cells.Item(13,19).Value2 = 0.0006770850438405203
$ws.Cells.Item(13,20).Value2 = 0.0006770850438405204

# Row 14
$ws.Cells.Item(14,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(14,2).Value2 = "Tgfb2"
$ws.Cells.Item(14,3).Value2 = "Tgfbr3"
$ws.Cells.Item(14,4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(14,5).Value2 = 1
$ws.Cells.Item(14,6).Value2 = 0.3333333333333333
$ws.Cells.Item(14,7).Value2 = 0.03423166666666667
$ws.Cells.Item(14,8).Value2 = 0.102695
$ws.Cells.Item(14,9).Value2 = 0.001475017749550216
$ws.Cells.Item(14,10).Value2 = 0.001475017749550216
$ws.Cells.Item(14,11).Value2 = 3
$ws.Cells.Item(14,12).Value2 = 1
$ws.Cells.Item(14,13).Value2 = 0.2401453333333333
$ws.Cells.Item(14,14).Value2 = 0.720436
$ws.Cells.Item(14,15).Value2 = 0.002698185430636104
$ws.Cells.Item(14,16).Value2 = 0.002698185430636104
$ws.Cells.Item(14,17).Value2 = 0.008220575002222221
$ws.Cells.Item(14,18).Value2 = 0.07398517501999999
$ws.Cells.Item(14,19).Value2 = 0.000003979871401766047
$ws.Cells.Item(14,20).Value2 = 0.000003979871401766048

# Row 15
$ws.Cells.Item(15,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(15,2).Value2 = "Tgfb2"
$ws.Cells.Item(15,3).Value2 = "Tgfbr3"
$ws.Cells.Item(15,4).Value2 = "MuSCs"
$ws.Cells.Item(15,5).Value2 = 1
$ws.Cells.Item(15,6).Value2 = 0.3333333333333333
$ws.Cells.Item(15,7).Value2 = 0.03423166666666667
$ws.Cells.Item(15,8).Value2 = 0.102695
$ws.Cells.Item(15,9).Value2 = 0.001475017749550216
$ws.Cells.Item(15,10).Value2 = 0.001475017749550216
$ws.Cells.Item(15,11).Value2 = 3
$ws.Cells.Item(15,12).Value2 = 1
$ws.Cells.Item(15,13).Value2 = 6.288874666666666
$ws.Cells.Item(15,14).Value2 = 18.866624
$ws.Cells.Item(15,15).Value2 = 0.0706595034147231
$ws.Cells.Item(15,16).Value2 = 0.0706595034147231
$ws.Cells.Item(15,17).Value2 = 0.2152786612977778
$ws.Cells.Item(15,18).Value2 = 1.93750795168
$ws.Cells.Item(15,19).Value2 = 0.0001042240217111207
$ws.Cells.Item(15,20).Value2 = 0.0001042240217111207

# Row 16
$ws.Cells.Item(16,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(16,2).Value2 = "Tgfb2"
$ws.Cells.Item(16,3).Value2 = "Tgfbr3"
$ws.Cells.Item(16,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(16,5).Value2 = 1
$ws.Cells.Item(16,6).Value2 = 0.3333333333333333
$ws.Cells.Item(16,7).Value2 = 0.03423166666666667
$ws.Cells.Item(16,8).Value2 = 0.102695
$ws.Cells.Item(16,9).Value2 = 0.001475017749550216
$ws.Cells.Item(16,10).Value2 = 0.001475017749550216
$ws.Cells.Item(16,11).Value2 = 1
$ws.Cells.Item(16,12).Value2 = 0.3333333333333333
$ws.Cells.Item(16,13).Value2 = 0.010449
$ws.Cells.Item(16,14).Value2 = 0.031347
$ws.Cells.Item(16,15).Value2 = 0.0001174011552645203
$ws.Cells.Item(16,16).Value2 = 0.0001174011552645203
$ws.Cells.Item(16,17).Value2 = 0.000357686685
$ws.Cells.Item(16,18).Value2 = 0.003219180165
$ws.Cells.Item(16,19).Value2 = 0.0000001731687878328683
$ws.Cells.Item(16,20).Value2 = 0.0000001731687878328683

# Row 17
$ws.Cells.Item(17,1).Value2 = "MuSCs"
$ws.Cells.Item(17,2).Value2 = "Tgfb2"
$ws.Cells.Item(17,3).Value2 = "Tgfbr3"
$ws.Cells.Item(17,4).Value2 = "ECs"
$ws.Cells.Item(17,5).Value2 = 3
$ws.Cells.Item(17,6).Value2 = 1
$ws.Cells.Item(17,7).Value2 = 7.520146
$ws.Cells.Item(17,8).Value2 = 22.560438
$ws.Cells.Item(17,9).Value2 = 0.324037650203293
$ws.Cells.Item(17,10).Value2 = 0.3240376502032931
$ws.Cells.Item(17,11).Value2 = 3
$ws.Cells.Item(17,12).Value2 = 1
$ws.Cells.Item(17,13).Value2 = 41.607769
$ws.Cells.Item(17,14).Value2 = 124.823307
$ws.Cells.Item(17,15).Value2 = 0.4674897261536314
$ws.Cells.Item(17,16).Value2 = 0.4674897261536314
$ws.Cells.Item(17,17).Value2 = 312.896497614274
$ws.Cells.Item(17,18).Value2 = 2816.068478528466
$ws.Cells.Item(17,19).Value2 = 0.1514842723570037
$ws.Cells.Item(17,20).Value2 = 0.1514842723570037

# Row 18
$ws.Cells.Item(18,1).Value2 = "MuSCs"
$ws.Cells.Item(18,2).Value2 = "Tgfb2"
$ws.Cells.Item(18,3).Value2 = "Tgfbr3"
$ws.Cells.Item(18,4).Value2 = "FAPs"
$ws.Cells.Item(18,5).Value2 = 3
$ws.Cells.Item(18,6).Value2 = 1
$ws.Cells.Item(18,7).Value2 = 7.520146
$ws.Cells.Item(18,8).Value2 = 22.560438
$ws.Cells.Item(18,9).Value2 = 0.324037650203293
$ws.Cells.Item(18,10).Value2 = 0.3240376502032931
$ws.Cells.Item(18,11).Value2 = 3
$ws.Cells.Item(18,12).Value2 = 1
$ws.Cells.Item(18,13).Value2 = 40.85529333333333
$ws.Cells.Item(18,14).Value2 = 122.56588
$ws.Cells.Item(18,15).Value2 = 0.4590351838457449
$ws.Cells.Item(18,16).Value2 = 0.4590351838457449
$ws.Cells.Item(18,17).Value2 = 307.2377707394933
$ws.Cells.Item(18,18).Value2 = 2765.13993665544
$ws.Cells.Item(18,19).Value2 = 0.1487446823340118
$ws.Cells.Item(18,20).Value2 = 0.1487446823340118

# Row 19
$ws.Cells.Item(19,1).Value2 = "MuSCs"
$ws.Cells.Item(19,2).Value2 = "Tgfb2"
$ws.Cells.Item(19,3).Value2 = "Tgfbr3"
$ws.Cells.Item(19,4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(19,5).Value2 = 3
$ws.Cells.Item(19,6).Value2 = 1
$ws.Cells.Item(19,7).Value2 = 7.520146
$ws.Cells.Item(19,8).Value2 = 22.560438
$ws.Cells.Item(19,9).Value2 = 0.324037650203293
$ws.Cells.Item(19,10).Value2 = 0.3240376502032931
$ws.Cells.Item(19,11).Value2 = 3
$ws.Cells.Item(19,12).Value2 = 1
$ws.Cells.Item(19,13).Value2 = 0.2401453333333333
$ws.Cells.Item(19,14).Value2 = 0.720436
$ws.Cells.Item(19,15).Value2 = 0.002698185430636104
$ws.Cells.Item(19,16).Value2 = 0.002698185430636104
$ws.Cells.Item(19,17).Value2 = 1.805927967885333
$ws.Cells.Item(19,18).Value2 = 16.253351710968
$ws.Cells.Item(19,19).Value2 = 0.0008743136667560834
$ws.Cells.Item(19,20).Value2 = 0.0008743136667560835

# Row 20
$ws.Cells.Item(20,1).Value2 = "MuSCs"
$ws.Cells.Item(20,2).Value2 = "Tgfb2"
$ws.Cells.Item(20,3).Value2 = "Tgfbr3"
$ws.Cells.Item(20,4).Value2 = "MuSCs"
$ws.Cells.Item(20,5).Value2 = 3
$ws.Cells.Item(20,6).Value2 = 1
$ws.Cells.Item(20,7).Value2 = 7.520146
$ws.Cells.Item(20,8).Value2 = 22.560438
$ws.Cells.Item(20,9).Value2 = 0.324037650203293
$ws.Cells.Item(20,10).Value2 = 0.3240376502032931
$ws.Cells.Item(20,11).Value2 = 3
$ws.Cells.Item(20,12).Value2 = 1
$ws.Cells.Item(20,13).Value2 = 6.288874666666666
$ws.Cells.Item(20,14).Value2 = 18.866624
$ws.Cells.Item(20,15).Value2 = 0.0706595034147231
$ws.Cells.Item(20,16).Value2 = 0.0706595034147231
$ws.Cells.Item(20,17).Value2 = 47.29325566903466
$ws.Cells.Item(20,18).Value2 = 425.6393010213119
$ws.Cells.Item(20,19).Value2 = 0.02289633945103843
$ws.Cells.Item(20,20).Value2 = 0.02289633945103844

# Row 21
$ws.Cells.Item(21,1).Value2 = "MuSCs"
$ws.Cells.Item(21,2).Value2 = "Tgfb2"
$ws.Cells.Item(21,3).Value2 = "Tgfbr3"
$ws.Cells.Item(21,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(21,5).Value2 = 3
$ws.Cells.Item(21,6).Value2 = 1
$ws.Cells.Item(21,7).Value2 = 7.520146
$ws.Cells.Item(21,8).Value2 = 22.560438
$ws.Cells.Item(21,9).Value2 = 0.324037650203293
$ws.Cells.Item(21,10).Value2 = 0.3240376502032931
$ws.Cells.Item(21,11).Value2 = 1
$ws.Cells.Item(21,12).Value2 = 0.3333333333333333
$ws.Cells.Item(21,13).Value2 = 0.010449
$ws.Cells.Item(21,14).Value2 = 0.031347
$ws.Cells.Item(21,15).Value2 = 0.0001174011552645203
$ws.Cells.Item(21,16).Value2 = 0.0001174011552645203
$ws.Cells.Item(21,17).Value2 = 0.07857800555399999
$ws.Cells.Item(21,18).Value2 = 0.7072020499859999
$ws.Cells.Item(21,19).Value2 = 0.00003804239448306712
$ws.Cells.Item(21,20).Value2 = 0.00003804239448306713

